$d = $word.ActiveDocument

$replacements = @(
    @{old="371÷5="; new="258÷2="},
    @{old="846÷4="; new="475÷9="},
    @{old="135÷2="; new="427÷7="},
    @{old="690÷7="; new="129÷3="},
    @{old="996÷7="; new="854÷9="},
    @{old="389÷5="; new="329÷3="},
    @{old="469÷6="; new="212÷4="},
    @{old="862÷2="; new="753÷4="},
    @{old="225÷6="; new="222÷6="},
    @{old="322÷4="; new="127÷7="},
    @{old="639÷9="; new="545÷8="},
    @{old="732÷5="; new="963÷3="},
    @{old="931÷5="; new="736÷2="},
    @{old="107÷5="; new="666÷5="},
    @{old="205÷3="; new="654÷5="},
    @{old="876÷3="; new="983÷9="},
    @{old="419÷4="; new="898÷9="},
    @{old="103÷3="; new="449÷9="},
    @{old="148÷3="; new="470÷8="},
    @{old="635÷8="; new="133÷9="},
    @{old="273÷8="; new="988÷6="},
    @{old="804÷9="; new="315÷5="},
    @{old="174÷5="; new="298÷5="},
    @{old="658÷5="; new="183÷7="},
    @{old="269÷4="; new="317÷2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
